# "list of external ontologies" -- fill in the full ontology names for the
# second half of the list (rows 9-22), which previously only had the
# acronym in column A and an empty column B.
#
# Shared-string order matters for a byte-faithful round trip: every new
# B-column value is appended to the shared-string table in the order it is
# first assigned, and the original commit's sharedStrings.xml shows "Human
# Phenotype Ontology" landing at the very end of the table (index 43)
# even though its row (13, acronym HP) sits in the middle of the block.
# So we set B13 last, after all the other new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value  = "Experimental Factor Ontology"
$ws.Range("B10").Value = "Environment Ontology"
$ws.Range("B11").Value = "Physico-Chemical Methods and Properties"
$ws.Range("B12").Value = "Gene Ontology"
$ws.Range("B14").Value = "Information Artifact Ontology"
$ws.Range("B15").Value = "National Cancer Innstitute Thesaurus"
$ws.Range("B16").Value = "NanoParticle Ontology"
$ws.Range("B17").Value = "Ontology of Adverse Events"
$ws.Range("B18").Value = "Ontology of Biological and Clinical Statistics"
$ws.Range("B19").Value = "Ontology for Biomedical Investigation"
$ws.Range("B20").Value = "Phenotype quality Ontology"
$ws.Range("B21").Value = "Statistics Ontology"
$ws.Range("B22").Value = "Unit of Measurement Ontology"
$ws.Range("B13").Value = "Human Phenotype Ontology"

# Column B widened very slightly (cosmetic re-save) and a narrow column D
# stub (defaultColWidth) appears in the saved file.
$ws.Columns.Item(2).ColumnWidth = 52.5
$ws.Columns.Item(4).ColumnWidth = 8

# Selection now spans the whole table instead of resting on B9.
[void]$ws.Range("A1:B22").Select()

# Page setup gains explicit paper size / orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
